$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 2 4 '29.177.69'
Set-TextValue 2 5 '  +0.85%  '
Set-TextValue 3 4 '1.833.54'
Set-TextValue 3 5 '  -0.03%  '
Set-TextValue 4 4 '0.9996'
Set-TextValue 4 5 '  +0.04%  '
Set-TextValue 5 4 '240.37'
Set-TextValue 5 5 '  -1.98%  '
Set-TextValue 6 4 '0.6849'
Set-TextValue 6 5 '  -1.13%  '
Set-TextValue 8 4 '0.3015'
Set-TextValue 8 5 '  -1.00%  '
Set-TextValue 9 4 '0.07469'
Set-TextValue 9 5 '  -2.82%  '
Set-TextValue 10 4 '23.13'
Set-TextValue 10 5 '  -1.03%  '
Set-TextValue 11 4 '0.07659'
Set-TextValue 11 5 '  -1.92%  '
Set-TextValue 12 4 '1.831.06'
Set-TextValue 12 5 '  -0.11%  '
Set-TextValue 13 4 '5.060'
Set-TextValue 13 5 '  -0.86%  '
Set-TextValue 14 4 '0.6825'
Set-TextValue 14 5 '  +0.08%  '
Set-TextValue 15 4 '86.98'
Set-TextValue 15 5 '  -6.68%  '
Set-TextValue 16 4 '6.213'
Set-TextValue 16 5 '  -5.56%  '
Set-TextValue 17 4 '29.160.61'
Set-TextValue 17 5 '  +0.79%  '
Set-TextValue 18 4 '0.000008187'
Set-TextValue 18 5 '  -0.97%  '
Set-TextValue 19 4 '2.082.86'
Set-TextValue 19 5 '  +0.41%  '
Set-TextValue 20 4 '12.54'
Set-TextValue 20 5 '  -1.21%  '
Set-TextValue 21 4 '226.55'
Set-TextValue 21 5 '  -5.96%  '
Set-TextValue 22 5 '  +0.07%  '
Set-TextValue 23 4 '7.414'
Set-TextValue 23 5 '  -0.62%  '
Set-TextValue 24 4 '1.000'
Set-TextValue 24 5 '  +0.05%  '
Set-TextValue 25 4 '0.1459'
Set-TextValue 25 5 '  -2.73%  '
Set-TextValue 26 4 '159.61'
Set-TextValue 26 5 '  +0.50%  '
Set-TextValue 27 4 '8.771'
Set-TextValue 27 5 '  +0.13%  '
Set-TextValue 28 4 '18.07'
Set-TextValue 28 5 '  -0.59%  '
Set-TextValue 29 4 '1.507'
Set-TextValue 29 5 '  -2.27%  '
Set-TextValue 30 5 '  +1.06%  '
Set-TextValue 31 5 '  -0.19%  '
Set-TextValue 32 4 '1.205'
Set-TextValue 32 5 '  +0.98%  '
Set-TextValue 33 4 '0.05150'
Set-TextValue 33 5 '  +0.81%  '
Set-TextValue 34 4 '0.7675'
Set-TextValue 34 5 '  -1.48%  '
Set-TextValue 35 4 '1.847'
Set-TextValue 35 5 '  -0.54%  '
Set-TextValue 36 4 '1.132'
Set-TextValue 36 5 '  -1.18%  '
Set-TextValue 37 5 '  -0.86%  '
Set-TextValue 38 4 '1.308.30'
Set-TextValue 38 5 '  +1.34%  '
Set-TextValue 39 4 '0.01836'
Set-TextValue 39 5 '  -1.23%  '
Set-TextValue 40 4 '2.710'
Set-TextValue 40 5 '  +0.45%  '
Set-TextValue 41 4 '0.9346'
Set-TextValue 41 5 '  -2.07%  '
Set-TextValue 42 4 '5.828'
Set-TextValue 42 5 '  -5.55%  '
Set-TextValue 43 4 '104.32'
Set-TextValue 43 5 '  -2.32%  '
Set-TextValue 44 4 '1.000'
Set-TextValue 44 5 '  +0.09%  '
Set-TextValue 45 4 '65.57'
Set-TextValue 45 5 '  +2.50%  '
Set-TextValue 46 2 'EnergySwap'
Set-TextValue 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 46 4 '9.609'
Set-TextValue 46 5 '  -0.77%  '
Set-TextValue 47 2 'RocketPoolETH'
Set-TextValue 47 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 47 4 '1.982.21'
Set-TextValue 47 5 '  +0.36%  '
Set-TextValue 48 2 'Mantle'
Set-TextValue 48 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 48 4 '0.5200'
Set-TextValue 48 5 '  +0.52%  '
Set-TextValue 49 4 '1.771'
Set-TextValue 49 5 '  +0.84%  '
Set-TextValue 50 4 '0.00000000121'
Set-TextValue 50 5 '  -1.74%  '
Set-TextValue 51 4 '0.07337'
Set-TextValue 51 5 '  +20.49%  '
